# Applies the commit "Added story, end of room and you died narrations":
#   - drops the stray _GoBack bookmark from the "Take5" line (Word relocates
#     it into the middle of the new "Story lines:" paragraph below)
#   - appends the new Story lines / End of room lines / You died sections that
#     follow the existing "Intro Lines:" block
$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Re-write the "Take5" paragraph without its bookmark -------------
$take5Range = $d.Content
$take5Range.Find.Execute("Take5 – Constipated Nic Cage", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$take5Range.Expand(4) | Out-Null   # wdParagraph -> include the paragraph mark
$take5Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Take5 – Constipated Nic Cage</w:t></w:r></w:p>')

# --- 2. Build the new paragraphs to append -------------------------------
$introGap = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$storyHeading = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Story lines:</w:t></w:r></w:p>'
$storyPara = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">While Kevin was having the slumber of a lifetime, his friends Alex Dungeon Crusher and Jonny the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Brauny</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> were </w:t></w:r><w:r><w:t>chopping down</w:t></w:r><w:r><w:t xml:space="preserve"> mons</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">ters. </w:t></w:r><w:r><w:t xml:space="preserve">However, they left the striders because they didn’t want eye goop on their swords. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Kevin, not having anything else to do, decided he would finish the job. Yeah… well, we’ll see how that goes.</w:t></w:r></w:p>'
$storyGap = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$endOfRoomHeading = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">End of room lines:</w:t></w:r></w:p>'
$endOfRoomLine1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Hm I’m surprised he killed that one</w:t></w:r></w:p>'
$endOfRoomLine2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Lucky shot I guess</w:t></w:r></w:p>'
$endOfRoomLine3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I think that one was already weakened</w:t></w:r></w:p>'
$endOfRoomLine4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I guess Striders really aren’t that strong</w:t></w:r></w:p>'
$youDiedGap = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$youDiedHeading = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">You died:</w:t></w:r></w:p>'
$youDiedLine1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I guess Kevin should have stayed in bed today</w:t></w:r></w:p>'
$youDiedLine2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I guess this is a job for a real hero</w:t></w:r></w:p>'
$youDiedLine3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To be honest, I didn’t think he’d make it this far</w:t></w:r></w:p>'
$youDiedLine4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Yeah, well he tried his best</w:t></w:r></w:p>'

$newContentXml = @($introGap, $storyHeading, $storyPara, $storyGap, $endOfRoomHeading, $endOfRoomLine1, $endOfRoomLine2, $endOfRoomLine3, $endOfRoomLine4, $youDiedGap, $youDiedHeading, $youDiedLine1, $youDiedLine2, $youDiedLine3, $youDiedLine4) -join ''

# --- 3. Append the new content at the very end of the document body ------
$endRange = $d.Content
$endRange.Collapse(0) | Out-Null   # wdCollapseEnd
$endRange.InsertXML($newContentXml)
